# Updated cryptos list (GitHub Actions data refresh).
# Updates Price (column D) and Volume(1h) (column E) figures for the coin
# rows in Sheet1, and swaps the BinanceUSD / Stellar rows (26 <-> 27),
# per the latest data pull.
#
# D-column values are prefixed with a leading apostrophe so Excel keeps
# them as literal text (matching the original "General"-formatted,
# text-typed cells) instead of re-interpreting strings like "215.40" or
# "1.00" as numbers and dropping trailing zeros / the thousands dots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'25.944.53"
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').Value = "'1.638.15"
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = "'215.40"
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('D9').Value = "'0.0639"
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').Value = "'19.63"
$ws.Range('E10').Value = '  -1.90%  '
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = "'1.864.41"
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('D14').Value = "'1.639.80"
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('E15').Value = '  -1.06%  '
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = "'62.98"
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').Value = "'25.947.52"
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').Value = "'193.02"
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('D21').Value = "'4.36"
$ws.Range('E21').Value = '  -1.60%  '
$ws.Range('E22').Value = '  -1.30%  '
$ws.Range('D23').Value = "'6.29"
$ws.Range('D24').Value = "'144.15"
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('E25').Value = '  +1.14%  '

# Row 26/27 swap: BinanceUSD and Stellar traded places in the ranking.
$ws.Range('B26').Value = 'BinanceUSD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = "'0.130"
$ws.Range('E27').Value = '  +3.65%  '
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('D29').Value = "'15.54"
$ws.Range('E29').Value = '  -0.51%  '
$ws.Range('D30').Value = "'1.24"
$ws.Range('E30').Value = '  -0.60%  '
$ws.Range('D31').Value = "'0.0503"
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('E34').Value = '  -4.15%  '
$ws.Range('D35').Value = "'2.45"
$ws.Range('E35').Value = '  +1.78%  '
$ws.Range('D36').Value = "'0.902"
$ws.Range('E36').Value = '  -1.26%  '
$ws.Range('D37').Value = "'1.138.28"
$ws.Range('E37').Value = '  +0.46%  '
$ws.Range('E38').Value = '  -1.60%  '
$ws.Range('D39').Value = "'2.47"
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('E40').Value = '  +0.31%  '
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('E42').Value = '  -3.04%  '
$ws.Range('D43').Value = "'99.30"
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').Value = "'1.774.29"
$ws.Range('E45').Value = '  -0.47%  '
$ws.Range('E46').Value = '  +2.21%  '
$ws.Range('E47').Value = '  +0.54%  '
$ws.Range('E48').Value = '  +3.12%  '
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('D50').Value = "'7.69"
$ws.Range('D51').Value = "'0.415"
$ws.Range('E51').Value = '  -0.71%  '